$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add the 7 new "MS River" gauge-station rows (23-29).
# Column A = short station code, Column B = long station name,
# Column C = a numeric reading (all 0 for now, same as the other new rows).
# ---------------------------------------------------------------------------

$stationCodes = @("WestPoint", "Alliance", "Carrollton", "BCSpillway", "BCSpillwayN", "Reserve")
$stationNames = @(
    "MS River at West Point a la Hache",
    "MS River at Alliance",
    "MS River at Carrollton",
    "MS River at Bonnet Carre Spillway",
    "MS River at Bonnet Carre Spillway North",
    "MS River at Reserve"
)

# Column A values first (so the shared-string table picks up the six station
# codes before the six descriptions - mirrors how the source workbook was
# built), then column B.
for ($i = 0; $i -lt $stationCodes.Length; $i++) {
    $row = 23 + $i
    $ws.Range("A$row").Value = $stationCodes[$i]
}
for ($i = 0; $i -lt $stationNames.Length; $i++) {
    $row = 23 + $i
    $ws.Range("B$row").Value = $stationNames[$i]
}

# The 7th new station (Morgan City) is entered as a standalone A/B pair.
$ws.Range("A29").Value = "MorganCity"
$ws.Range("B29").Value = "Lower Atchafalaya River at Morgan City"

# Column C numeric readings for every new row - same numeric format (one
# decimal place) as the existing data rows above, copied from C22.
foreach ($row in 23..29) {
    $ws.Range("C$row").Value = 0
}
$ws.Range("C22").Copy() | Out-Null
$ws.Range("C23:C29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Style: column A for the new rows gets a dedicated font (Arial 10, black)
# and an 8-decimal numeric format, matching the style added for this import.
# Build it once on A23, then copy the resulting formatting (not value) to
# the other new rows via copy / paste-special so the same style (cellXf)
# entry is reused instead of creating a new one per cell.
# ---------------------------------------------------------------------------

$ws.Range("A23").Font.Name = "Arial"
$ws.Range("A23").Font.Size = 10
$ws.Range("A23").Font.Color = 0
$ws.Range("A23").NumberFormat = "0.00000000"

$ws.Range("A23").Copy() | Out-Null
$ws.Range("A24:A29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Selection moves to the newly-added last row, like it would right after
# typing it in.
# ---------------------------------------------------------------------------
$ws.Range("A29:C29").Select()
